$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44313
$ws.Range("K4").Value = 'Winter Nelis'
$ws.Range("L4").Value = 'Tercera'
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15500
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("S4").Value = 861

# Row 5
$ws.Range("K5").Value = 'Packham''s Triumph'
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 19000
$ws.Range("P5").Value = 18500
$ws.Range("Q5").Value = '$/caja 18 kilos empedrada'
$ws.Range("S5").Value = 1028

# Row 6
$ws.Range("D6").Value = 44474
$ws.Range("K6").Value = 'Winter Nelis'
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("S6").Value = 972

# Row 7
$ws.Range("D7").Value = 44280
$ws.Range("M7").Value = 350
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 24500
$ws.Range("S7").Value = 1361

# Row 8
$ws.Range("D8").Value = 44280
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 24500
$ws.Range("S8").Value = 1361

# Row 9
$ws.Range("D9").Value = 44329
$ws.Range("M9").Value = 340
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("S9").Value = 1194

# Row 10
$ws.Range("D10").Value = 44769
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 16500
$ws.Range("S10").Value = 917

# Row 11
$ws.Range("D11").Value = 44769
$ws.Range("K11").Value = 'Winter Nelis'
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 15500
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("S11").Value = 861

# Row 12
$ws.Range("D12").Value = 44525
$ws.Range("K12").Value = 'Packham''s Triumph'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("S12").Value = 1083

# Row 13
$ws.Range("D13").Value = 44642
$ws.Range("M13").Value = 270
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("S13").Value = 1083

# Row 14
$ws.Range("D14").Value = 44292
$ws.Range("K14").Value = 'Packham''s Triumph'
$ws.Range("N14").Value = 22000
$ws.Range("O14").Value = 23000
$ws.Range("P14").Value = 22500
$ws.Range("S14").Value = 1250

# Row 15
$ws.Range("D15").Value = 44292
$ws.Range("K15").Value = 'Winter Nelis'
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 22000
$ws.Range("O15").Value = 23000
$ws.Range("P15").Value = 22500
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("S15").Value = 1250

# Row 16
$ws.Range("D16").Value = 44323
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 16000
$ws.Range("P16").Value = 15500
$ws.Range("Q16").Value = '$/bandeja 18 kilos granel'
$ws.Range("S16").Value = 861

# Row 18
$ws.Range("D18").Value = 44341
$ws.Range("K18").Value = 'Packham''s Triumph'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 17000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 17500
$ws.Range("Q18").Value = '$/caja 18 kilos granel'
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 972

# Row 19
$ws.Range("D19").Value = 44747
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 270
$ws.Range("N19").Value = 19000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 19500
$ws.Range("Q19").Value = '$/bandeja 18 kilos granel'
$ws.Range("S19").Value = 1083

# Row 20
$ws.Range("D20").Value = 44747
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 19000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19500
$ws.Range("Q20").Value = '$/bandeja 18 kilos granel'
$ws.Range("S20").Value = 1083

# Row 21
$ws.Range("D21").Value = 44715
$ws.Range("R21").Value = 'Región de O''Higgins'

# Row 22
$ws.Range("D22").Value = 44715
$ws.Range("K22").Value = 'Winter Nelis'
$ws.Range("N22").Value = 17000
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 17500
$ws.Range("Q22").Value = '$/caja 18 kilos granel'
$ws.Range("S22").Value = 972

# Row 23
$ws.Range("D23").Value = 44678
$ws.Range("K23").Value = 'Packham''s Triumph'
$ws.Range("N23").Value = 17000
$ws.Range("O23").Value = 18000
$ws.Range("P23").Value = 17500
$ws.Range("Q23").Value = '$/caja 18 kilos granel'
$ws.Range("S23").Value = 972

# Row 24
$ws.Range("D24").Value = 44763
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 300
$ws.Range("P24").Value = 17500
$ws.Range("Q24").Value = '$/caja 18 kilos granel'
$ws.Range("S24").Value = 972

# Row 25
$ws.Range("D25").Value = 44763
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 300
$ws.Range("Q25").Value = '$/caja 18 kilos granel'

# Row 26
$ws.Range("D26").Value = 44371
$ws.Range("L26").Value = 'Calibre 90'
$ws.Range("M26").Value = 140
$ws.Range("N26").Value = 17000
$ws.Range("O26").Value = 18000
$ws.Range("P26").Value = 17429
$ws.Range("Q26").Value = '$/caja 18 kilos embalada'
$ws.Range("S26").Value = 968

# Row 27
$ws.Range("D27").Value = 44371
$ws.Range("L27").Value = 'Calibre 80'
$ws.Range("M27").Value = 120
$ws.Range("N27").Value = 17000
$ws.Range("O27").Value = 18000
$ws.Range("P27").Value = 17500
$ws.Range("Q27").Value = '$/caja 18 kilos embalada'
$ws.Range("S27").Value = 972

# Row 28
$ws.Range("D28").Value = 44355
$ws.Range("K28").Value = 'Packham''s Triumph'
$ws.Range("M28").Value = 200
$ws.Range("N28").Value = 17000
$ws.Range("O28").Value = 18000
$ws.Range("P28").Value = 17500
$ws.Range("Q28").Value = '$/caja 18 kilos granel'
$ws.Range("R28").Value = 'Región Metropolitana'
$ws.Range("S28").Value = 972

# Row 29
$ws.Range("D29").Value = 44355
$ws.Range("K29").Value = 'Winter Nelis'
$ws.Range("M29").Value = 250
$ws.Range("N29").Value = 17000
$ws.Range("O29").Value = 18000
$ws.Range("P29").Value = 17500
$ws.Range("Q29").Value = '$/caja 18 kilos granel'
$ws.Range("R29").Value = 'Región Metropolitana'
$ws.Range("S29").Value = 972

# Row 30
$ws.Range("D30").Value = 44336
$ws.Range("M30").Value = 250
$ws.Range("N30").Value = 21000
$ws.Range("O30").Value = 22000
$ws.Range("P30").Value = 21500
$ws.Range("S30").Value = 1194

# Row 31
$ws.Range("D31").Value = 44421
$ws.Range("M31").Value = 270
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 17000
$ws.Range("P31").Value = 16500
$ws.Range("S31").Value = 917

# Row 32
$ws.Range("D32").Value = 44421
$ws.Range("N32").Value = 16000
$ws.Range("O32").Value = 17000
$ws.Range("P32").Value = 16500
$ws.Range("S32").Value = 917

# Row 33
$ws.Range("D33").Value = 44699
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 300
$ws.Range("N33").Value = 17000
$ws.Range("O33").Value = 18000
$ws.Range("P33").Value = 17500
$ws.Range("Q33").Value = '$/caja 18 kilos empedrada'
$ws.Range("S33").Value = 972

# Row 34
$ws.Range("D34").Value = 44699
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 250
$ws.Range("N34").Value = 17000
$ws.Range("O34").Value = 18000
$ws.Range("P34").Value = 17500
$ws.Range("Q34").Value = '$/caja 18 kilos empedrada'
$ws.Range("S34").Value = 972

# Row 35
$ws.Range("D35").Value = 44497
$ws.Range("M35").Value = 300
$ws.Range("Q35").Value = '$/bandeja 18 kilos granel'
$ws.Range("R35").Value = 'Región de O''Higgins'

# Row 36
$ws.Range("D36").Value = 44497
$ws.Range("Q36").Value = '$/bandeja 18 kilos granel'
$ws.Range("R36").Value = 'Región de O''Higgins'
